# The edit permutes the D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado) and P (Precio $/Kg)
# values across data rows 2..44 - i.e. each row ends up showing the
# date/volume/price figures that used to belong to a different row,
# while every other column (Mercado, Region, Categoria, Unidad, Origen,
# Kg o Unidades, Clasificacion, ...) stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destinationRow -> sourceRow (source row's original values move into
# the destination row)
$rowMap = @{
    2  = 12
    3  = 41
    4  = 16
    5  = 26
    6  = 21
    7  = 39
    8  = 6
    9  = 15
    10 = 42
    11 = 23
    12 = 8
    13 = 44
    14 = 25
    15 = 2
    16 = 22
    17 = 36
    18 = 27
    19 = 33
    20 = 20
    21 = 5
    22 = 17
    23 = 14
    24 = 29
    25 = 28
    26 = 38
    27 = 30
    28 = 31
    29 = 7
    30 = 9
    31 = 34
    32 = 32
    33 = 19
    34 = 18
    35 = 3
    36 = 40
    37 = 10
    38 = 35
    39 = 4
    40 = 37
    41 = 11
    42 = 24
    43 = 43
    44 = 13
}

$cols = @("D", "J", "K", "L", "M", "P")

# Snapshot every original value first - this is a permutation, so we
# must read all the "before" values before writing any "after" ones.
# Value2() is used (rather than Value) so the raw date serial number
# comes back for column D instead of a formatted date string.
$original = @{}
foreach ($col in $cols) {
    for ($row = 2; $row -le 44; $row++) {
        $original[$col + $row] = $ws.Range($col + $row).Value2()
    }
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    foreach ($col in $cols) {
        $ws.Range($col + $destRow).Value = $original[$col + $srcRow]
    }
}
